$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PipelineState")
$ws.Range("A3").Value = "SkyBox-Cubed"
$ws.Range("B3").Value = "SkyBox-Cubed"
